# Add a new "2022" data column (column S) to the disasters-deaths table,
# mirroring the formatting of the existing "2021" column (R) for every row,
# and filling in the 2022 figures reported in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> 2022 value (numbers, or "-" where MES KR reported no data).
$updates = @(
    @{Row=3;  Value=$null},
    @{Row=4;  Value=2022},
    @{Row=5;  Value=135},
    @{Row=6;  Value=99},
    @{Row=7;  Value=36},
    @{Row=8;  Value=97},
    @{Row=9;  Value=80},
    @{Row=10; Value=17},
    @{Row=11; Value=17},
    @{Row=12; Value=11},
    @{Row=13; Value=6},
    @{Row=14; Value=5},
    @{Row=15; Value=3},
    @{Row=16; Value=2},
    @{Row=17; Value='-'},
    @{Row=18; Value='-'},
    @{Row=19; Value='-'},
    @{Row=20; Value=6},
    @{Row=21; Value=1},
    @{Row=22; Value=5},
    @{Row=23; Value='-'},
    @{Row=24; Value='-'},
    @{Row=25; Value='-'},
    @{Row=26; Value=10},
    @{Row=27; Value=4},
    @{Row=28; Value=6},
    @{Row=29; Value='-'},
    @{Row=30; Value='-'},
    @{Row=31; Value='-'},
    @{Row=32; Value='-'},
    @{Row=33; Value='-'},
    @{Row=34; Value='-'}
)

foreach ($u in $updates) {
    $r = $u.Row
    $src = $ws.Cells.Item($r, 18)   # column R (18th column) - the 2021 column
    $dst = $ws.Cells.Item($r, 19)   # column S (19th column) - new 2022 column

    # Copy R's formatting (font/border/number-format) onto S first.
    $src.Copy()
    $dst.PasteSpecial(-4122)

    if ($null -ne $u.Value) {
        $dst.Value = $u.Value
    }
}

# Match the saved selection from the authored workbook.
$ws.Range("T4").Select()

Write-Output "2022 column added"
